$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "58.480.32"
$ws.Cells.Item(2,5).Value = "  -3.98%  "

# Row 3
$ws.Cells.Item(3,4).Value = "2.645.91"
$ws.Cells.Item(3,5).Value = "  -2.12%  "

# Row 4
$ws.Cells.Item(4,5).Value = "  +0.06%  "

# Row 5
$ws.Cells.Item(5,4).Value = "'521.34"
$ws.Cells.Item(5,5).Value = "  -1.00%  "

# Row 6
$ws.Cells.Item(6,4).Value = "'144.20"
$ws.Cells.Item(6,5).Value = "  -0.60%  "

# Row 7
$ws.Cells.Item(7,4).Value = "'0.999"
$ws.Cells.Item(7,5).Value = "  +0.27%  "

# Row 8
$ws.Cells.Item(8,5).Value = "  -1.59%  "

# Row 9
$ws.Cells.Item(9,4).Value = "'6.70"
$ws.Cells.Item(9,5).Value = "  -0.04%  "

# Row 10
$ws.Cells.Item(10,5).Value = "  -3.11%  "

# Row 11
$ws.Cells.Item(11,5).Value = "  -0.86%  "

# Row 12
$ws.Cells.Item(12,5).Value = "  +1.48%  "

# Row 13
$ws.Cells.Item(13,4).Value = "3.106.65"
$ws.Cells.Item(13,5).Value = "  -2.36%  "

# Row 14
$ws.Cells.Item(14,4).Value = "58.445.33"
$ws.Cells.Item(14,5).Value = "  -3.95%  "

# Row 15
$ws.Cells.Item(15,5).Value = "  -2.45%  "

# Row 16
$ws.Cells.Item(16,5).Value = "  -1.45%  "

# Row 17
$ws.Cells.Item(17,4).Value = "2.643.96"
$ws.Cells.Item(17,5).Value = "  -7.56%  "

# Row 18
$ws.Cells.Item(18,4).Value = "'338.03"
$ws.Cells.Item(18,5).Value = "  -2.89%  "

# Row 19
$ws.Cells.Item(19,5).Value = "  -2.25%  "

# Row 20
$ws.Cells.Item(20,5).Value = "  -1.38%  "

# Row 21
$ws.Cells.Item(21,4).Value = "'6.30"
$ws.Cells.Item(21,5).Value = "  -1.66%  "

# Row 22
$ws.Cells.Item(22,5).Value = "  +0.32%  "

# Row 23
$ws.Cells.Item(23,4).Value = "'64.68"

# Row 24
$ws.Cells.Item(24,4).Value = "'0.423"
$ws.Cells.Item(24,5).Value = "  +0.65%  "

# Row 25
$ws.Cells.Item(25,5).Value = "  -2.12%  "

# Row 26
$ws.Cells.Item(26,4).Value = "'0.997"
$ws.Cells.Item(26,5).Value = "  +0.51%  "

# Row 27
$ws.Cells.Item(27,4).Value = "0.0₃0798"
$ws.Cells.Item(27,5).Value = "  -2.61%  "

# Row 28
$ws.Cells.Item(28,5).Value = "  -2.75%  "

# Row 29
$ws.Cells.Item(29,4).Value = "'6.62"
$ws.Cells.Item(29,5).Value = "  -2.12%  "

# Row 30
$ws.Cells.Item(30,5).Value = "  +0.03%  "

# Row 31
$ws.Cells.Item(31,5).Value = "  -1.39%  "

# Row 32
$ws.Cells.Item(32,4).Value = "'152.61"
$ws.Cells.Item(32,5).Value = "  +1.54%  "

# Row 33
$ws.Cells.Item(33,5).Value = "  -1.51%  "

# Row 34
$ws.Cells.Item(34,5).Value = "  -2.71%  "

# Row 35
$ws.Cells.Item(35,2).Value = "ImmutableX"
$ws.Cells.Item(35,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35,4).Value = "'1.18"
$ws.Cells.Item(35,5).Value = "  -4.24%  "

# Row 36
$ws.Cells.Item(36,2).Value = "SuiNetwork"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(36,4).Value = "'0.903"
$ws.Cells.Item(36,5).Value = "  -4.70%  "

# Row 37
$ws.Cells.Item(37,4).Value = "'0.856"
$ws.Cells.Item(37,5).Value = "  -2.59%  "

# Row 38
$ws.Cells.Item(38,4).Value = "'36.71"
$ws.Cells.Item(38,5).Value = "  -0.77%  "

# Row 39
$ws.Cells.Item(39,4).Value = "'1.45"
$ws.Cells.Item(39,5).Value = "  -4.97%  "

# Row 40
$ws.Cells.Item(40,5).Value = "  -1.01%  "

# Row 41
$ws.Cells.Item(41,5).Value = "  +0.27%  "

# Row 42
$ws.Cells.Item(42,4).Value = "'0.607"
$ws.Cells.Item(42,5).Value = "  -0.96%  "

# Row 43
$ws.Cells.Item(43,4).Value = "'0.0967"
$ws.Cells.Item(43,5).Value = "  -2.40%  "

# Row 44
$ws.Cells.Item(44,4).Value = "'268.87"
$ws.Cells.Item(44,5).Value = "  -5.62%  "

# Row 45
$ws.Cells.Item(45,4).Value = "'19.35"
$ws.Cells.Item(45,5).Value = "  -4.16%  "

# Row 46
$ws.Cells.Item(46,5).Value = "  +1.62%  "

# Row 47
$ws.Cells.Item(47,5).Value = "  -0.58%  "

# Row 48
$ws.Cells.Item(48,4).Value = "2.039.69"
$ws.Cells.Item(48,5).Value = "  -4.99%  "

# Row 49
$ws.Cells.Item(49,2).Value = "RenderToken"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49,4).Value = "'4.70"
$ws.Cells.Item(49,5).Value = "  -2.25%  "

# Row 50
$ws.Cells.Item(50,2).Value = "VeChain"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50,4).Value = "'0.0228"
$ws.Cells.Item(50,5).Value = "  -3.16%  "

# Row 51
$ws.Cells.Item(51,4).Value = "'18.31"
$ws.Cells.Item(51,5).Value = "  -4.63%  "
